$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Populate column B (Progress) with values for rows 2-14 (one trial row per row)
$progress = @(0.25, 0.5, 0.875, 0.75, 0, 0.125, 0.25, 0, 0, 0.125, 0.25, 0, 0)
for ($i = 0; $i -lt $progress.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $progress[$i]
}

# Update the selection to match the new state
$ws.Range("A1:B32").Select()
